$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.302.99'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '2.290.12'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '155.82'
$ws.Range("E5").Value = '  +15,465.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.32'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '96.86'
$ws.Range("E7").Value = '  +4.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.498'
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '35.65'
$ws.Range("E11").Value = '  +8.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0808'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '2.644.37'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.58'
$ws.Range("E16").Value = '  +1.83%  '
$ws.Range("D17").Value = '2.281.21'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.798'
$ws.Range("E18").Value = '  +4.31%  '
$ws.Range("D19").Value = '42.216.41'
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.81'
$ws.Range("E20").Value = '  +4.11%  '
$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.31'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '244.58'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.62'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.25'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.35'
$ws.Range("E29").Value = '  +6.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.74'
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.21'
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.38'
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0756'
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.11'
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("E37").Value = '  +4.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.38'
$ws.Range("E38").Value = '  +2.85%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.23'
$ws.Range("E42").Value = '  +7.54%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.01'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.014.72'
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("E45").Value = '  +11.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0285'
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.27'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.01'
$ws.Range("E48").Value = '  +2.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.93'
$ws.Range("E49").Value = '  +3.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.54'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.88'
$ws.Range("E51").Value = '  -0.16%  '
